{"js": "// Spanish localization pass for the \"Partner email \u2013 reminder to submit\n// documents\" template. The same English phrase always maps to the same\n// Spanish phrase *except* for two short connector fragments (\", at \" and\n// \" or \") which also appear, unchanged, inside a later paragraph that the\n// source diff leaves untranslated. Those two are therefore replaced only\n// inside the specific paragraphs that the diff actually touches; every\n// other phrase is replaced document-wide (all occurrences).\n\n// Simple whole-document text replacements: every occurrence of `from`\n// becomes `to` (safe because distinct source strings never map to more\n// than one target string in this diff).\nconst wholeDocReplacements = [\n  [\"English\", \"Ingl\u00e9s\"],\n  [\" / Portuguese / French / Thai / Vietnamese / Spanish\", \" / Portugu\u00e9s / Franc\u00e9s / Tailand\u00e9s / Vietnamita / Espa\u00f1ol\"],\n  [\"Brief\", \"Breve\"],\n  [\"An email sent to partners in the target country who RSVPed yes but haven\\u2019t sent their documents to us. It will be sent via customer.io\",\n   \"An email sent to partners in the target country who RSVPed yes but haven\\u2019t sent their documents to us. Se enviar\u00e1 a trav\u00e9s de customer.io\"],\n  [\"Target audience\", \"P\u00fablico objetivo\"],\n  [\"Subject line\", \"Asunto\"],\n  [\" \\u2014 have you submitted your docs?  \", \" - \u00bfhas presentado tus documentos?  \"],\n  [\"Don\\u2019t forget to send your documents\", \"No olvides enviar tus documentos\"],\n  [\"Hi \", \"Hola \"],\n  [\"We\\u2019re excited to see you at the upcoming \", \"Estamos emocionados por verte en el pr\u00f3ximo \"],\n  [\"To confirm your registration, we need the following documents from you by \",\n   \"Para confirmar tu inscripci\u00f3n, necesitamos que nos env\u00edes los siguientes documentos antes del \"],\n  [\"Please send a copy of these documents to your country manager, \",\n   \"Por favor, env\u00eda una copia de estos documentos al responsable de tu pa\u00eds, \"],\n  [\" (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.\",\n   \" (WhatsApp), para que podamos realizar los preparativos necesarios para ti, incluidos el alojamiento y el transporte.\"],\n  [\"If you have any questions, please contact your country manager.\",\n   \"Si tienes alguna pregunta, entra en contacto con el responsable de tu pa\u00eds.\"],\n  [\"We look forward to seeing you there!\", \"\u00a1Esperamos verte ah\u00ed!\"],\n  [\"Dear \", \"Estimado \"],\n  [\"To ensure you have the best experience at this event, we need the following documents from you by \",\n   \"Para asegurarnos de que disfrutes de la mejor experiencia en este evento, necesitamos que nos env\u00edes los siguientes documentos hasta el \"],\n  [\"Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.\",\n   \"Por favor, responde a este correo electr\u00f3nico con una copia de estos documentos para que podamos hacer los arreglos necesarios para ti, incluyendo alojamiento y transporte.\"],\n  [\"If you have any questions, please contact us via \", \"Si tienes alguna pregunta, entra en contacto con nosotros por \"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of wholeDocReplacements) {\n  const results = body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Paragraph-scoped replacements for the two ambiguous fragments: the diff\n// only translates the first \"..., at [EMAIL ADDRESS] or [WHATSAPP NO]...\"\n// block (the customer.io partner email) and the \"...live chat or\n// WhatsApp...\" sentence; the later, near-duplicate \"country manager, [NAME],\n// at [EMAIL ADDRESS] or [WHATSAPP NO]\" paragraph is left in English.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nasync function replaceInParagraph(paragraphText, from, to, occurrenceIndex) {\n  paragraphText.load(\"text\");\n  await context.sync();\n  const results = paragraphText.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > (occurrenceIndex || 0)) {\n    results.items[occurrenceIndex || 0].insertText(to, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text.indexOf(\"Please send a copy of these documents to your country manager\") !== -1 ||\n      text.indexOf(\"Por favor, env\\u00eda una copia de estos documentos al responsable de tu pa\\u00eds\") !== -1) {\n    // This is the partner-email paragraph that still has English \", at \" / \" or \".\n    await replaceInParagraph(para, \", at \", \", a \", 0);\n    await replaceInParagraph(para, \" or \", \" o \", 0);\n  } else if (text.indexOf(\"please contact us via\") !== -1 ||\n             text.indexOf(\"entra en contacto con nosotros por\") !== -1) {\n    // The \"...live chat or WhatsApp...\" sentence.\n    await replaceInParagraph(para, \" or \", \" o \", 0);\n  }\n  // NOTE: the later paragraph \"If you have any questions, please contact\n  // your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO]\n  // (WhatsApp).\" is intentionally left untouched \u2014 the source diff does\n  // not translate it.\n}\n\n// The comment text (\"choose either one\" -> \"elija uno de los dos\").\nconst comments = body.getComments();\ncomments.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < comments.items.length; i++) {\n  comments.items[i].load(\"content\");\n}\nawait context.sync();\nfor (let i = 0; i < comments.items.length; i++) {\n  if (comments.items[i].content === \"choose either one\") {\n    comments.items[i].content = \"elija uno de los dos\";\n  }\n}\nawait context.sync();\n", "ps1": "# Spanish localization pass for the \"Partner email \u2013 reminder to submit\n# documents\" template.\n#\n# Strategy: every English phrase that needs translating maps to exactly one\n# Spanish phrase everywhere it occurs in the body -- EXCEPT for two short\n# connector fragments (\", at \" and \" or \") that also occur, unchanged,\n# inside a later paragraph which the source diff leaves in English. Those\n# two fragments are therefore replaced only within the specific paragraphs\n# the diff actually touches (using a paragraph-scoped Range.Find), while\n# everything else is replaced document-wide with Replace:=wdReplaceAll.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null\n}\n\nfunction Replace-TextInParagraph($paraIndex, $findText, $replaceText) {\n    $p = $d.Paragraphs.Item($paraIndex)\n    $rng = $p.Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null\n}\n\n# --- Whole-document replacements (safe: same source text always maps to\n#     the same target text everywhere it appears) ---\nReplace-AllText \"English\" \"Ingl\u00e9s\"\nReplace-AllText \" / Portuguese / French / Thai / Vietnamese / Spanish\" \" / Portugu\u00e9s / Franc\u00e9s / Tailand\u00e9s / Vietnamita / Espa\u00f1ol\"\nReplace-AllText \"Brief\" \"Breve\"\nReplace-AllText \"An email sent to partners in the target country who RSVPed yes but haven't sent their documents to us. It will be sent via customer.io\" \"An email sent to partners in the target country who RSVPed yes but haven't sent their documents to us. Se enviar\u00e1 a trav\u00e9s de customer.io\"\nReplace-AllText \"Target audience\" \"P\u00fablico objetivo\"\nReplace-AllText \"Subject line\" \"Asunto\"\nReplace-AllText \" \u2014 have you submitted your docs?  \" \" - \u00bfhas presentado tus documentos?  \"\nReplace-AllText \"Don't forget to send your documents\" \"No olvides enviar tus documentos\"\nReplace-AllText \"Hi \" \"Hola \"\nReplace-AllText \"We're excited to see you at the upcoming \" \"Estamos emocionados por verte en el pr\u00f3ximo \"\nReplace-AllText \"To confirm your registration, we need the following documents from you by \" \"Para confirmar tu inscripci\u00f3n, necesitamos que nos env\u00edes los siguientes documentos antes del \"\nReplace-AllText \"Please send a copy of these documents to your country manager, \" \"Por favor, env\u00eda una copia de estos documentos al responsable de tu pa\u00eds, \"\nReplace-AllText \" (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.\" \" (WhatsApp), para que podamos realizar los preparativos necesarios para ti, incluidos el alojamiento y el transporte.\"\nReplace-AllText \"If you have any questions, please contact your country manager.\" \"Si tienes alguna pregunta, entra en contacto con el responsable de tu pa\u00eds.\"\nReplace-AllText \"We look forward to seeing you there!\" \"\u00a1Esperamos verte ah\u00ed!\"\nReplace-AllText \"Dear \" \"Estimado \"\nReplace-AllText \"To ensure you have the best experience at this event, we need the following documents from you by \" \"Para asegurarnos de que disfrutes de la mejor experiencia en este evento, necesitamos que nos env\u00edes los siguientes documentos hasta el \"\nReplace-AllText \"Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.\" \"Por favor, responde a este correo electr\u00f3nico con una copia de estos documentos para que podamos hacer los arreglos necesarios para ti, incluyendo alojamiento y transporte.\"\nReplace-AllText \"If you have any questions, please contact us via \" \"Si tienes alguna pregunta, entra en contacto con nosotros por \"\n\n# --- Paragraph-scoped replacements ---\n# Find the two paragraphs by distinctive surrounding text rather than a\n# hard-coded index, since that is robust to the prior replacements above\n# having already run.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Por favor, env\u00eda una copia de estos documentos al responsable de tu pa\u00eds*\") {\n        # \"..., at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp), para que...\" (customer.io partner email)\n        Replace-TextInParagraph $i \", at \" \", a \"\n        Replace-TextInParagraph $i \" or \" \" o \"\n    }\n    elseif ($t -like \"*entra en contacto con nosotros por*\") {\n        # \"...live chat or WhatsApp...\"\n        Replace-TextInParagraph $i \" or \" \" o \"\n    }\n    # NOTE: the later paragraph \"If you have any questions, please contact\n    # your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO]\n    # (WhatsApp).\" is intentionally left untouched \u2014 the source diff does\n    # not translate it.\n}\n\n# --- Comment text ---\nfor ($i = 1; $i -le $d.Comments.Count; $i++) {\n    $c = $d.Comments.Item($i)\n    if ($c.Range.Text -eq \"choose either one\") {\n        $c.Range.Text = \"elija uno de los dos\"\n    }\n}\n"}
